$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# --- Insert two new rows into the "Tableau1" structured table ---
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Copy the formatting of the first data row (row 2) down onto the two new rows (3,4)
$ws.Range("A2:G2").Copy() | Out-Null
$ws.Range("A3:G4").PasteSpecial(-4122) | Out-Null

# --- Fill in the new row values, column by column (matches authoring order) ---
$ws.Range("A3").Value = "M"
$ws.Range("A4").Value = "Mlle"

$ws.Range("B3").Value = "Lezou Marman"
$ws.Range("B4").Value = "Kouassi Linda"

$ws.Range("C3").Value = "Agronomix"
$ws.Range("C4").Value = "Agronomix"

$ws.Range("D3").Value = "10/12/1996"
$ws.Range("D4").Value = "10/13/1990"

$ws.Range("E3").Value = "lm@bbmlm.ci"
$ws.Range("E4").Value = "kl@bbmlm.ci"

$ws.Range("F3").Value = "NOVUS CUSTOMER"
$ws.Range("F4").Value = "NOVUS CUSTOMER"

$ws.Range("G3").Value = "ZDdkNDNkMDJjNjBhZDg1NTc5YThiNT"
$ws.Range("G4").Value = "ZDdkNDNkMDJjNjBhZDg1NTc5YThiNT"

# --- Data validation lists must now cover rows 2 through 4 ---
$ws.Range("A2:A4").Validation.Delete()
$ws.Range("F2:F4").Validation.Delete()
$ws.Range("C2:C4").Validation.Delete()

$ws.Range("A2:A4").Validation.Add(3, 1, 1, "civility")
$ws.Range("F2:F4").Validation.Add(3, 1, 1, "nature")
$ws.Range("C2:C4").Validation.Add(3, 1, 1, "corporation")

# --- Hyperlinks: drop the old one on E2 and attach fresh ones to the new email cells ---
$ws.Range("E2").Hyperlinks.Delete()
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").VerticalAlignment = -4108

$ws.Range("E3").WrapText = $false
$ws.Range("E4").WrapText = $false

$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:lm@bbmlm.ci") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:kl@bbmlm.ci") | Out-Null

# --- Update the active selection left behind on the sheet ---
$ws.Activate() | Out-Null
$ws.Range("D10").Select() | Out-Null
